# Apply the 2025 adherent list changes:
#  - Insert a new "Régulier" column after "Statut"/"Date dernière adhésion" column,
#    right before "Envoi mails" (this shifts every following column right by one,
#    which matches the rest of the diff automatically).
#  - Rename / shorten several headers.
#  - Fix "Prenom" -> "Prénom" typo.
#  - Resize several columns.
#  - Set the sample row's new "Régulier" cell to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new column before column G ("Envoi mails").
$ws.Columns.Item(7).Insert()

# 2) Fix the "Prenom" -> "Prénom" header.
$ws.Range("C1").Value = "Prénom"

# 3) Rename headers to their new, shorter labels (post-insert layout).
$ws.Range("F1").Value = "Dern. adh."
$ws.Range("G1").Value = "Régulier"
$ws.Range("M1").Value = "Adh. 2024"
$ws.Range("N1").Value = "Adh. 2025"
$ws.Range("O1").Value = "Adh. 2026"
$ws.Range("P1").Value = "Dons 2025"
$ws.Range("Q1").Value = "Total 2025"
$ws.Range("R1").Value = "Dern. paiement"

# 4) Set the sample data row's new "Régulier" value (not a regular donor).
$ws.Range("G2").Value = $false

# 5) Resize the columns to match the new layout.
$widths = @{
    1  = 25
    2  = 15
    3  = 15
    4  = 10
    5  = 5
    6  = 8
    7  = 8
    8  = 8
    9  = 25
    10 = 8
    11 = 25
    12 = 10
    13 = 8
    14 = 8
    15 = 8
    16 = 8
    17 = 8
    18 = 10
    19 = 10
    20 = 20
}

foreach ($colIndex in $widths.Keys) {
    $ws.Columns.Item($colIndex).ColumnWidth = $widths[$colIndex] - 0.8333333333333334
}
